# products added to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "WalMArt"
$ws.Range("C1").Value = "Dollar Trap"
$ws.Range("D1").Value = "Office Repo"
$ws.Range("A1:D1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# ---- Product rows (A2:A16) ----
$products = @(
    "Ball Point Pen",
    "TI-35 Calculator",
    "100 Page Notebook",
    "8 oz Glue",
    "Clear Tape",
    "Eraser",
    "10 No.2 Pencils",
    "2inch Binder",
    "USB Stick 5gb",
    "Color Markers",
    "Stapler",
    "Planner Book",
    "Protractor",
    "Compass",
    "Liquid Paper"
)

for ($i = 0; $i -lt $products.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $products[$i]
}

# ---- Prices (B2:D16) ----
$prices = @(
    @(0.5, 0.4, 1.4),
    @(28, 33, 31),
    @(1.8, 1, 2),
    @(1.2, 0.8, 1.5),
    @(2.4, 1.4, 2.4),
    @(0.9, 0.2, 0.8),
    @(0.99, 0.59, 2.59),
    @(1.25, 3.25, 2.15),
    @(9.5, 14, 13),
    @(4.55, 2.55, 6),
    @(4.2, 2.2, 3),
    @(3.9, 5, 8),
    @(1, 2, 1),
    @(1.75, 2, 1),
    @(2, 1, 3)
)

for ($i = 0; $i -lt $prices.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $prices[$i][0]
    $ws.Cells.Item($row, 3).Value = $prices[$i][1]
    $ws.Cells.Item($row, 4).Value = $prices[$i][2]
}

# ---- Currency style for the price grid ----
$ws.Range("B2:D16").Style = "Currency"

# Keep the Currency style's font matching the workbook's default font
# (Excel reuses the normal 11pt font for this built-in style).
$curStyle = $wb.Styles.Item("Currency")
$curStyle.Font.Size = 11

$ws.Range("G11").Select() | Out-Null
